# EIA Table 1.16.A — roll the report forward from "October 2016/2015" to
# "November 2016/2015" (2017-01-31 update), and refresh the November data
# for Mountain / Nevada / Utah / Wyoming / Pacific Contiguous / California /
# Oregon / Pacific Noncontiguous / Hawaii / U.S. Total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write literal text into a cell without Excel's autoconvert turning
# month/year-looking strings (e.g. "November 2016") into date serials.
function Set-TextValue($range, [string]$text) {
    $origFmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $origFmt
}

# --- Title / subtitle strings ---------------------------------------------
# (plain sentence text; no risk of Excel's date autoconvert, so a direct
# assignment is enough here)

$ws.Range("A2").Value = "by State, by Sector, November 2016 and 2015 (Thousand Megawatthours)"

# Column headers on row 6 (period labels), every "October 2016"/"October 2015"
# occurrence across the All Sectors / Electric Power / Commercial / Industrial
# sub-columns.
foreach ($col in @("B","E","G","I","K")) {
    Set-TextValue $ws.Range($col + "6") "November 2016"
}
foreach ($col in @("C","F","H","J","L")) {
    Set-TextValue $ws.Range($col + "6") "November 2015"
}

# --- Data updates -----------------------------------------------------------

# Row 52: Mountain
$ws.Range("B52").Value = 397
$ws.Range("C52").Value = 362
$ws.Range("D52").Value = 0.096
$ws.Range("E52").Value = 23
$ws.Range("F52").Value = 21
$ws.Range("G52").Value = 374
$ws.Range("H52").Value = 341

# Row 55: Idaho
$ws.Range("D55").Value = 0.143

# Row 57: Nevada
$ws.Range("B57").Value = 339
$ws.Range("C57").Value = 319
$ws.Range("D57").Value = 0.063
$ws.Range("G57").Value = 339
$ws.Range("H57").Value = 319

# Row 59: Utah
$ws.Range("B59").Value = 48
$ws.Range("C59").Value = 35
$ws.Range("D59").Value = 0.382
$ws.Range("E59").Value = 23
$ws.Range("F59").Value = 21
$ws.Range("G59").Value = 25
$ws.Range("H59").Value = 14

# Row 61: Pacific Contiguous
$ws.Range("B61").Value = 1084
$ws.Range("C61").Value = 954
$ws.Range("D61").Value = 0.137
$ws.Range("E61").Value = 72
$ws.Range("G61").Value = 1013
$ws.Range("H61").Value = 881

# Row 62: California
$ws.Range("B62").Value = 1064
$ws.Range("C62").Value = 935
$ws.Range("D62").Value = 0.138
$ws.Range("E62").Value = 70
$ws.Range("F62").Value = 72
$ws.Range("G62").Value = 994
$ws.Range("H62").Value = 863

# Row 63: Oregon
$ws.Range("B63").Value = 20
$ws.Range("C63").Value = 19
$ws.Range("D63").Value = 0.064
$ws.Range("G63").Value = 18
$ws.Range("H63").Value = 18

# Row 65: Pacific Noncontiguous
$ws.Range("B65").Value = 26
$ws.Range("D65").Value = 0.449
$ws.Range("G65").Value = 26

# Row 67: Hawaii
$ws.Range("B67").Value = 26
$ws.Range("D67").Value = 0.449
$ws.Range("G67").Value = 26

# Row 68: U.S. Total
$ws.Range("B68").Value = 1507
$ws.Range("C68").Value = 1334
$ws.Range("D68").Value = 0.13
$ws.Range("E68").Value = 94
$ws.Range("F68").Value = 93
$ws.Range("G68").Value = 1412
$ws.Range("H68").Value = 1240
